$wb = $excel.ActiveWorkbook

# Update the "deletePatientInfo" sheet: change A2 value and selection,
# and make it the active/selected tab.
$ws2 = $wb.Worksheets.Item("deletePatientInfo")
$ws2.Range("A2").Value = 94
$ws2.Activate()
$ws2.Range("B2").Select()
